# Testexercise5.xlsx - "calculate_price_test" sheet update
#
# 1) Row 6 (test case 5 - "Earlier than easter holiday, second tier"):
#    the purchase date (column G) moved from 1/1/2022 to 10/4/2021.
# 2) A new test case (row 22, ID 21) is appended: "0 on first day" /
#    "Between 9-18 April,  total=0", total=0, date 9/4/2022, expected=0.
# 3) The sheet's saved view (selection) is updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix the date of purchase for test case 5 -------------------------
$ws.Range("G6").Value = 44296

# --- 2) Append new row 22 -------------------------------------------------
# Row 2 carries the same cell-style pattern (s=3/5/5/4/4/4/12/4/4/4/4) that
# the new row needs, so copy its formatting only, then fill in the values.
$ws.Range("A2:K2").Copy()
$ws.Range("A22").PasteSpecial(-4122)

$ws.Range("A22").Value = 21
$ws.Range("B22").Value = "0 on first day"
$ws.Range("C22").Value = "Between 9-18 April,  total=0"
$ws.Range("D22").Value = "Exercise5/vEP1_1"
$ws.Range("E22").Value = "order exist"
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 44660
$ws.Range("H22").Value = 0

# --- 3) Update the recorded selection / scroll position -------------------
$ws.Range("A10:XFD10").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
